# Append new scrape results (2025-08-27 12:36 JST) to the top of the job list
# and record a new stats row, mirroring the scraper's usual "prepend latest
# run, keep history below" behaviour.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ランサーズ" (job listing) - insert two new rows at the top of the
# data (row 2), pushing the existing 9 rows down to rows 4-12.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Rows("2:3").Insert()

# Row-insert shifts cell values/styles correctly but leaves the Hyperlinks
# collection pointing at the old (pre-insert) addresses, so drop them all and
# rebuild from scratch once every cell is in its final place.
$ws1.Hyperlinks.Delete()

# --- New row 2 ---
$ws1.Range("A2").Value = "2025-08-27 12:36:04"
$ws1.Range("B2").Value = "【急募】outsystemsでホテル予約アプリの開発依頼"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Range("F2").Value = "https://www.lancers.jp/work/detail/5380830"
$ws1.Range("G2").Value = 85
$ws1.Range("H2").Value = "◆開発 ◇アプリ"

# --- New row 3 ---
$ws1.Range("A3").Value = "2025-08-27 12:36:04"
$ws1.Range("B3").Value = "wordpressレンダリングを妨げるリソースの除外"
$ws1.Range("C3").Value = "システム開発"
$ws1.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws1.Range("E3").Value = "期限情報なし"
$ws1.Range("F3").Value = "https://www.lancers.jp/work/detail/5016989"
$ws1.Range("G3").Value = 33
$ws1.Range("H3").Value = "○WordPress"

# --- Rebuild the hyperlinks for every URL cell, rows 2-12 ---
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws1.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws1.Hyperlinks.Add($cell, $url)
}

# ---------------------------------------------------------------------------
# Sheet "統計" (stats) - append a new summary row for this run.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A3").Value = "2025-08-27T12:36:04.614974"
$ws2.Range("B3").Value = 11
$ws2.Range("C3").Value = "全案件リスト"
$ws2.Range("D3").Value = 72.7
$ws2.Range("E3").Value = 3
$ws2.Range("F3").Value = 6
$ws2.Range("G3").Value = 11
